# Update header row labels so Power BI can automatically treat the first
# row as a header (prefix the year labels with "Ano"/"Intervalo").

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5: "Fonte/Tecnologia" | 2015 | 2030 | 2040 | 2050 -> prefix "Ano "
$anoSheets = @(1, 2, 3, 5)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4: intervals -> prefix "Intervalo "
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 6: only has B1 = 2015 -> prefix "Ano "
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B1").Value = "Ano 2015"
